$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the defined names used by the new/updated formulas ---
$wb.Names.Add('Output', '=Sheet1!$B$4')
$wb.Names.Add('S_0', '=Sheet1!$D$5')
$wb.Names.Add('S_1', '=Sheet1!$D$6')
$wb.Names.Add('S_2', '=Sheet1!$D$7')
$wb.Names.Add('S_3', '=Sheet1!$D$8')
$wb.Names.Add('S0', '=Sheet1!$D$5')
$wb.Names.Add('tm', '=Sheet1!$B$3')

# --- Update the raw inputs in column B ---
$ws.Range("B6").Value = 780
$ws.Range("B7").Value = 800

# --- Rewrite the QUOTIENT chain in column D to use the new "tm" named range ---
$ws.Range("D6").Formula = '=QUOTIENT(D5*(tm-C6),tm)'
$ws.Range("D7").Formula = '=QUOTIENT(D6*(tm-C7),tm)'
$ws.Range("D8").Formula = '=QUOTIENT(D7*(tm-C8),tm)'

# --- New row 9: intermediate helper value ---
$ws.Range("B9").Formula = '=(tm-S_3-1)*Output'

# --- Row 10: new formula driven by the helper cell above ---
$ws.Range("B10").Formula = '=IF(B9>=0,QUOTIENT(B9,tm)+1,0)'

# --- Restore the selected cell as recorded in the saved workbook ---
$ws.Range("B7").Select() | Out-Null
